# Generate Report for Handoff
# Adds two new localization entries (d2525536-... and d77c1228-...) as new
# rows across the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$commitHash = "c013a34671323ee191ba6ee31bcadd44a62278ee"
$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/"

$file1 = "d2525536-afb4-49a1-9951-26382a7aad42.md"
$file2 = "d77c1228-e007-429b-a451-91b6e2b4bfc4.md"

$zhXlf1 = "d2525536-afb4-49a1-9951-26382a7aad42.3dfb4f9e625c28a8e2eddd68622cd89c604a4e98.zh-cn.xlf"
$zhXlf2 = "d77c1228-e007-429b-a451-91b6e2b4bfc4.e98d1301818359ff62cfd816ac2a994c2b3411c2.zh-cn.xlf"
$deXlf1 = "d2525536-afb4-49a1-9951-26382a7aad42.3dfb4f9e625c28a8e2eddd68622cd89c604a4e98.de-de.xlf"
$deXlf2 = "d77c1228-e007-429b-a451-91b6e2b4bfc4.e98d1301818359ff62cfd816ac2a994c2b3411c2.de-de.xlf"

$status = "Ready for handoff"
$overviewDate = "2016-09-05 10:21:59"
$zhHandoffDate = "2016-09-05 10:21:47"
$zeroDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Sheet "Overview" (table3 / displayName "Overview")
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item(1)

$tblOverview.ListRows.Add() | Out-Null
$wsOverview.Range("A4").Value = $file1
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), ($ghBase + $file1), "", "", ("e2e\" + $file1)) | Out-Null
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = $status
$wsOverview.Range("F4").Value = $status
$wsOverview.Range("G4").Value = $overviewDate
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$tblOverview.ListRows.Add() | Out-Null
$wsOverview.Range("A5").Value = $file2
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), ($ghBase + $file2), "", "", ("e2e\" + $file2)) | Out-Null
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = $status
$wsOverview.Range("F5").Value = $status
$wsOverview.Range("G5").Value = $overviewDate
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (table1 / displayName "zh_cn")
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$tblZh = $wsZh.ListObjects.Item(1)

$tblZh.ListRows.Add() | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), ($ghBase + $file1), "", "", $file1) | Out-Null
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = $zhXlf1
$wsZh.Range("H4").Value = $zhHandoffDate
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = ""
$wsZh.Range("J4").Value = ""
$wsZh.Range("K4").Value = $zeroDate
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = ""

$tblZh.ListRows.Add() | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), ($ghBase + $file2), "", "", $file2) | Out-Null
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = $status
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = $zhXlf2
$wsZh.Range("H5").Value = $zhHandoffDate
$wsZh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = $zeroDate
$wsZh.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "False"
$wsZh.Range("P5").Value = ""

# ---------------------------------------------------------------------------
# Sheet "de-de" (table2 / displayName "de_de")
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$tblDe = $wsDe.ListObjects.Item(1)

$tblDe.ListRows.Add() | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), ($ghBase + $file1), "", "", $file1) | Out-Null
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = $deXlf1
$wsDe.Range("H4").Value = $overviewDate
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = ""
$wsDe.Range("J4").Value = ""
$wsDe.Range("K4").Value = $zeroDate
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = ""

$tblDe.ListRows.Add() | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), ($ghBase + $file2), "", "", $file2) | Out-Null
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = $status
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = $deXlf2
$wsDe.Range("H5").Value = $overviewDate
$wsDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = $zeroDate
$wsDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "False"
$wsDe.Range("P5").Value = ""

Write-Host "Report generated for handoff."
